$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-21 (columns A-T).
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = numeric metric columns.
$csv = @"
ECs,Il16,Kcnj10,ECs,3,1,3.648214333333333,10.944643,0.1438182892240057,0.1438182892240057,1,0.3333333333333333,0.06268366666666667,0.188051,0.01538636255968474,0.01538636255968473,0.2286834511992222,2.058151060793,0.002212840340714152,0.002212840340714152
ECs,Il16,Kcnj10,Inflammatory-Mac,3,1,3.648214333333333,10.944643,0.1438182892240057,0.1438182892240057,3,1,1.727625,5.182875,0.4240636521556707,0.4240636521556707,6.302746287624999,56.724716588625,0.06098810897511241,0.06098810897511239
ECs,Il16,Kcnj10,MuSCs,3,1,3.648214333333333,10.944643,0.1438182892240057,0.1438182892240057,3,1,0.075253,0.225759,0.01847163708308845,0.01847163708308845,0.2745390732263333,2.470851659037,0.002656559244456484,0.002656559244456483
ECs,Il16,Kcnj10,Resolving-Mac,3,1,3.648214333333333,10.944643,0.1438182892240057,0.1438182892240057,3,1,2.208414,6.625242,0.5420783482015561,0.5420783482015561,8.056767608733999,72.510908478606,0.07796078066372267,0.07796078066372265
FAPs,Il16,Kcnj10,ECs,3,1,3.016808666666666,9.050426,0.1189272947567555,0.1189272947567555,1,0.3333333333333333,0.06268366666666667,0.188051,0.01538636255968474,0.01538636255968473,0.1891046288584444,1.701941659726,0.001829858475369934,0.001829858475369934
FAPs,Il16,Kcnj10,Inflammatory-Mac,3,1,3.016808666666666,9.050426,0.1189272947567555,0.1189272947567555,3,1,1.727625,5.182875,0.4240636521556707,0.4240636521556707,5.21191407275,46.90722665475,0.0504327429555437,0.0504327429555437
FAPs,Il16,Kcnj10,MuSCs,3,1,3.016808666666666,9.050426,0.1189272947567555,0.1189272947567555,3,1,0.075253,0.225759,0.01847163708308845,0.01847163708308845,0.2270239025926667,2.043215123334,0.002196781828020277,0.002196781828020276
FAPs,Il16,Kcnj10,Resolving-Mac,3,1,3.016808666666666,9.050426,0.1189272947567555,0.1189272947567555,3,1,2.208414,6.625242,0.5420783482015561,0.5420783482015561,6.662362494787999,59.961262453092,0.06446791149782162,0.06446791149782162
Inflammatory-Mac,Il16,Kcnj10,ECs,3,1,10.36785533333333,31.103566,0.4087169998040092,0.4087169998040092,1,0.3333333333333333,0.06268366666666667,0.188051,0.01538636255968474,0.01538636255968473,0.6498951877628889,5.849056689866,0.006288667943291081,0.00628866794329108
Inflammatory-Mac,Il16,Kcnj10,Inflammatory-Mac,3,1,10.36785533333333,31.103566,0.4087169998040092,0.4087169998040092,3,1,1.727625,5.182875,0.4240636521556707,0.4240636521556707,17.91176607025,161.20589463225,0.1733220236349967,0.1733220236349967
Inflammatory-Mac,Il16,Kcnj10,MuSCs,3,1,10.36785533333333,31.103566,0.4087169998040092,0.4087169998040092,3,1,0.075253,0.225759,0.01847163708308845,0.01847163708308845,0.7802122173993333,7.021909956594,0.007549672090068392,0.007549672090068391
Inflammatory-Mac,Il16,Kcnj10,Resolving-Mac,3,1,10.36785533333333,31.103566,0.4087169998040092,0.4087169998040092,3,1,2.208414,6.625242,0.5420783482015561,0.5420783482015561,22.896516868108,206.068651812972,0.221556636135653,0.221556636135653
MuSCs,Il16,Kcnj10,ECs,3,1,0.05357233333333333,0.160717,0.00211190479115806,0.00211190479115806,1,0.3333333333333333,0.06268366666666667,0.188051,0.01538636255968474,0.01538636255968473,0.003358110285222222,0.030222992567,3.24945328082932E-05,3.249453280829319E-05
MuSCs,Il16,Kcnj10,Inflammatory-Mac,3,1,0.05357233333333333,0.160717,0.00211190479115806,0.00211190479115806,3,1,1.727625,5.182875,0.4240636521556707,0.4240636521556707,0.092552902375,0.832976121375,0.0008955820587435462,0.0008955820587435462
MuSCs,Il16,Kcnj10,MuSCs,3,1,0.05357233333333333,0.160717,0.00211190479115806,0.00211190479115806,3,1,0.075253,0.225759,0.01847163708308845,0.01847163708308845,0.004031478800333333,0.036283309203,3.90103388563074E-05,3.901033885630739E-05
MuSCs,Il16,Kcnj10,Resolving-Mac,3,1,0.05357233333333333,0.160717,0.00211190479115806,0.00211190479115806,3,1,2.208414,6.625242,0.5420783482015561,0.5420783482015561,0.118309890946,1.064789018514,0.001144817860749914,0.001144817860749914
Resolving-Mac,Il16,Kcnj10,ECs,3,1,8.280381,24.841143,0.3264255114240716,0.3264255114240716,1,0.3333333333333333,0.06268366666666667,0.188051,0.01538636255968474,0.01538636255968473,0.519044642477,4.671401782293,0.005022501267501278,0.005022501267501277
Resolving-Mac,Il16,Kcnj10,Inflammatory-Mac,3,1,8.280381,24.841143,0.3264255114240716,0.3264255114240716,3,1,1.727625,5.182875,0.4240636521556707,0.4240636521556707,14.305393225125,128.748539026125,0.1384251945312744,0.1384251945312744
Resolving-Mac,Il16,Kcnj10,MuSCs,3,1,8.280381,24.841143,0.3264255114240716,0.3264255114240716,3,1,0.075253,0.225759,0.01847163708308845,0.01847163708308845,0.623123511393,5.608111602536999,0.006029613581686994,0.006029613581686992
Resolving-Mac,Il16,Kcnj10,Resolving-Mac,3,1,8.280381,24.841143,0.3264255114240716,0.3264255114240716,3,1,2.208414,6.625242,0.5420783482015561,0.5420783482015561,18.286509325734,164.578583931606,0.1769482020436089,0.1769482020436089
"@

$lines = $csv -split "`n"

$rowIndex = 2
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","

    # Columns A-D are text (cluster / gene symbol names)
    $ws.Cells.Item($rowIndex, 1).Value = $parts[0]
    $ws.Cells.Item($rowIndex, 2).Value = $parts[1]
    $ws.Cells.Item($rowIndex, 3).Value = $parts[2]
    $ws.Cells.Item($rowIndex, 4).Value = $parts[3]

    # Columns E-T are numeric (indices 4..19 in $parts map to columns 5..20)
    for ($col = 5; $col -le 20; $col++) {
        $strVal = $parts[$col - 1]
        $numVal = [double]$strVal
        $ws.Cells.Item($rowIndex, $col).Value = $numVal
    }

    $rowIndex++
}

Write-Output "Updated rows 2 through $($rowIndex - 1)"
